$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update quarterly data values per author's MV source refresh
# Row 97
$ws.Range("B97").Value = 2.5
$ws.Range("C97").Value = 2.9
$ws.Range("H97").Value = 1.2

# Row 98
$ws.Range("B98").Value = 1.3
$ws.Range("C98").Value = -1.4

# Row 99
$ws.Range("B99").Value = 2.5
$ws.Range("C99").Value = 5.2

# Row 100
$ws.Range("B100").Value = 2.8
$ws.Range("C100").Value = 4.7

# Row 101
$ws.Range("B101").Value = 2.6
$ws.Range("C101").Value = 1.8
$ws.Range("H101").Value = 1.8

# Row 102
$ws.Range("B102").Value = 3.8
$ws.Range("C102").Value = 3.3
$ws.Range("H102").Value = 2.7

# Row 103
$ws.Range("B103").Value = 3
$ws.Range("C103").Value = 2.3

# Row 104
$ws.Range("B104").Value = 2.2
$ws.Range("C104").Value = 1.3
$ws.Range("H104").Value = 1.6

# Row 105
$ws.Range("B105").Value = 1.9

# Row 106
$ws.Range("B106").Value = 1.6
$ws.Range("C106").Value = 2.4
$ws.Range("H106").Value = 2.2

# Row 107
$ws.Range("C107").Value = 1.2
$ws.Range("K107").Value = 4.2

# Row 108
$ws.Range("C108").Value = 2.4
$ws.Range("H108").Value = 1.6
$ws.Range("K108").Value = 4.5

# Row 109
$ws.Range("B109").Value = 2
$ws.Range("C109").Value = 2
$ws.Range("G109").Value = 2
$ws.Range("H109").Value = 3.2
$ws.Range("K109").Value = 4.8

# Row 110
$ws.Range("B110").Value = 1.9
$ws.Range("C110").Value = 1.9
$ws.Range("F110").Value = 3
$ws.Range("H110").Value = 2.7

# Row 111
$ws.Range("B111").Value = 2.1
$ws.Range("C111").Value = 2.3
$ws.Range("F111").Value = 1
$ws.Range("H111").Value = 3.2
$ws.Range("I111").Value = 2.8
$ws.Range("K111").Value = 5.7

# Row 112
$ws.Range("B112").Value = 2.3
$ws.Range("H112").Value = 2.7
$ws.Range("K112").Value = 6.1

# Row 113
$ws.Range("C113").Value = 3.8
$ws.Range("H113").Value = 3.6
$ws.Range("K113").Value = 5.8

# Row 114
$ws.Range("F114").Value = -0.1
$ws.Range("I114").Value = 2.9

# Row 115
$ws.Range("F115").Value = 0.2
$ws.Range("I115").Value = 3

# Row 117
$ws.Range("I117").Value = 3.3
$ws.Range("Q117").Value = -2.6

# Row 118
$ws.Range("F118").Value = 1.2
$ws.Range("I118").Value = 1.9
$ws.Range("Q118").Value = -2.5

# Row 119
$ws.Range("F119").Value = 0.8
$ws.Range("I119").Value = 2.3
$ws.Range("Q119").Value = -1.6

# Row 120
$ws.Range("I120").Value = 2.1

# Row 121
$ws.Range("F121").Value = -7.5
$ws.Range("I121").Value = 2.6
$ws.Range("Q121").Value = 6.1

# Row 122
$ws.Range("F122").Value = -2.3
$ws.Range("Q122").Value = 4.4

# Row 123
$ws.Range("F123").Value = -28.2
$ws.Range("I123").Value = -2.6
$ws.Range("Q123").Value = -10.4

# Row 124
$ws.Range("F124").Value = 22.8
$ws.Range("I124").Value = -1

# Row 125
$ws.Range("E125").Value = -1
$ws.Range("F125").Value = 11.9
$ws.Range("I125").Value = -1.1
$ws.Range("Q125").Value = 6.2

# Row 126
$ws.Range("E126").Value = -1.3
$ws.Range("F126").Value = -3.7
$ws.Range("M126").Value = -2.5
$ws.Range("P126").Value = 9.300000000000001
$ws.Range("Q126").Value = 7.2

# Row 127
$ws.Range("C127").Value = 6.6
$ws.Range("D127").Value = 22.2
$ws.Range("E127").Value = 7.5
$ws.Range("F127").Value = 1.3
$ws.Range("G127").Value = 13.6
$ws.Range("H127").Value = 8.199999999999999
$ws.Range("I127").Value = 6
$ws.Range("K127").Value = 16.1
$ws.Range("L127").Value = 11.1
$ws.Range("M127").Value = 8.199999999999999
$ws.Range("N127").Value = 10.3
$ws.Range("O127").Value = 7.5
$ws.Range("P127").Value = 7.4
$ws.Range("Q127").Value = 21.7
$ws.Range("S127").Value = 12.4
$ws.Range("T127").Value = 17.6
$ws.Range("U127").Value = 19.6
